$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: column B renamed from "Symbol" to "Building" ---
$ws.Range("B1").Value = "Building"

# --- Column B values (was empty, now holds a tower/block label) ---
# Rows 2-5 -> TF8, rows 6-9 -> TF9, rows 10-17 -> TF10
$ws.Range("B2").Value = "TF8"
$ws.Range("B3").Value = "TF8"
$ws.Range("B4").Value = "TF8"
$ws.Range("B5").Value = "TF8"

$ws.Range("B6").Value = "TF9"
$ws.Range("B7").Value = "TF9"
$ws.Range("B8").Value = "TF9"
$ws.Range("B9").Value = "TF9"

$ws.Range("B10").Value = "TF10"
$ws.Range("B11").Value = "TF10"
$ws.Range("B12").Value = "TF10"
$ws.Range("B13").Value = "TF10"
$ws.Range("B14").Value = "TF10"
$ws.Range("B15").Value = "TF10"
$ws.Range("B16").Value = "TF10"
$ws.Range("B17").Value = "TF10"

# --- B2:B5 need a distinct format: centered horizontally, no fill,
#     bordered, but WITHOUT vertical centering (unlike the rest of the
#     B column which keeps the existing row banding style). Clone the
#     plain bordered style from B6 (fontId 0, fillId 0, borderId 1) and
#     then strip the vertical alignment so only horizontal=center remains.
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B2:B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B2:B5").VerticalAlignment = -4107
$excel.CutCopyMode = 0

# --- Selection moves from H7 to B10:B17, active cell B10 ---
$ws.Range("B10:B17").Select() | Out-Null

# --- Phonetic info present on the sheet (noConversion) ---
$ws.SetPhonetic = $false
